$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EXAMPLESUB_DATALIST")

# Update the preprocessing pipeline folder/confounds-folder name
# from "preproc_fmriprep-1.5.8_MOD5" to "preproc_fmriprep-20.0.1"
# for rows 2 through 5, columns G (dataFolder) and H (confoundsFolder)
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 7).Value = "preproc_fmriprep-20.0.1"
    $ws.Cells.Item($r, 8).Value = "preproc_fmriprep-20.0.1"
}

# Update the active cell selection to match the new state
$ws.Range("H13").Select()
